$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing rows (old EXCL036/EXCL037 footnote rows) that are dropped entirely
$ws.Rows("58:59").Delete()

$ws.Range("C2").Value = 'INCL001'
$ws.Range("D2").Value = 'Inclusion Criteria'
$ws.Range("G2").Value = 'Signed Informed Consent Form'

$ws.Range("C3").Value = 'INCL002'
$ws.Range("D3").Value = 'Inclusion Criteria'
$ws.Range("G3").Value = 'Age ≥ 18 years at time of signing Informed Consent Form'

$ws.Range("C4").Value = 'INCL003'
$ws.Range("D4").Value = 'Inclusion Criteria'
$ws.Range("G4").Value = 'Ability to comply with the study protocol, including willingness to remain in the post-treatment period'

$ws.Range("C5").Value = 'INCL004'
$ws.Range("D5").Value = 'Inclusion Criteria'
$ws.Range("G5").Value = 'ECOG Performance Status of 0 or 1 (see Appendix 3)'

$ws.Range("C6").Value = 'INCL005'
$ws.Range("D6").Value = 'Inclusion Criteria'
$ws.Range("G6").Value = 'Histologically or cytologically documented NSCLC with unresectable Stage III NSCLC of either squamous or non-squamous histology Staging should be based on the 8th revised editio... (As per the protocol)'

$ws.Range("C7").Value = 'INCL006'
$ws.Range("D7").Value = 'Inclusion Criteria'
$ws.Range("G7").Value = 'Whole-body positron emission tomography (PET)-CT scan (from the base of skull to mid-thighs) for the purposes of staging, performed prior and within 42 days of the first dose of concurrent CRT'

$ws.Range("C8").Value = 'INCL007'
$ws.Range("D8").Value = 'Inclusion Criteria'
$ws.Range("G8").Value = 'At least two prior cycles of platinum-based chemotherapy concurrent with RT (CRT), which must be completed within 1 to 42 days prior to randomization in the study To ensure the... (As per the protocol)'

$ws.Range("C9").Value = 'INCL008'
$ws.Range("D9").Value = 'Inclusion Criteria'
$ws.Range("G9").Value = 'The RT component in the CRT must have been at a total dose of radiation of 60 (± 10) Gy (54−66 Gy) administered by IMRT (preferred) or 3D-conforming technique Sites are encourag... (As per the protocol)'

$ws.Range("C10").Value = 'INCL009'
$ws.Range("D10").Value = 'Inclusion Criteria'
$ws.Range("G10").Value = 'No progression during or following concurrent platinum-based CRT'

$ws.Range("C11").Value = 'INCL010'
$ws.Range("D11").Value = 'Inclusion Criteria'
$ws.Range("G11").Value = 'Tumor PD-L1 expression, as determined by SP263 IHC assay and documented by means of central testing of a representative tumor tissue, in either a previously obtained archival tu... (As per the protocol)'

$ws.Range("C12").Value = 'INCL011'
$ws.Range("D12").Value = 'Inclusion Criteria'
$ws.Range("G12").Value = 'Confirmed availability of representative tumor specimens in formalin- fixed, paraffin-embedded (FFPE) blocks (preferred) or at least 15−20 unstained serial slides, along with an... (As per the protocol)'

$ws.Range("C13").Value = 'INCL012'
$ws.Range("D13").Value = 'Inclusion Criteria'
$ws.Range("G13").Value = 'Life expectancy ≥ 12 weeks'

$ws.Range("C14").Value = 'INCL013'
$ws.Range("D14").Value = 'Inclusion Criteria'
$ws.Range("G14").Value = 'Adequate hematologic and end-organ function, defined by the following laboratory test results, obtained within 14 days prior to initiation of study treatment (Day 1 of Cycle 1):... (As per the protocol)'

$ws.Range("C15").Value = 'INCL014'
$ws.Range("D15").Value = 'Inclusion Criteria'
$ws.Range("G15").Value = 'Negative HIV test at screening'

$ws.Range("C16").Value = 'INCL015'
$ws.Range("D16").Value = 'Inclusion Criteria'
$ws.Range("G16").Value = 'Negative hepatitis B surface antigen (HBsAg) test at screening'

$ws.Range("C17").Value = 'INCL016'
$ws.Range("D17").Value = 'Inclusion Criteria'
$ws.Range("G17").Value = 'Positive hepatitis B surface antibody (HBsAb) test at screening, or negative HBsAb at screening accompanied by either of the following: – Negative hepatitis B core antibody (HBc... (As per the protocol)'

$ws.Range("C18").Value = 'INCL017'
$ws.Range("D18").Value = 'Inclusion Criteria'
$ws.Range("G18").Value = 'Negative hepatitis C virus (HCV) antibody test at screening, or positive HCV antibody test followed by a negative HCV RNA test at screening The HCV RNA test will be performed on... (As per the protocol)'

$ws.Range("C19").Value = 'INCL018'
$ws.Range("D19").Value = 'Inclusion Criteria'
$ws.Range("G19").Value = 'For women of childbearing potential: agreement to remain abstinent (refrain from heterosexual intercourse) or use contraception, and agreement to refrain from donating eggs, as... (As per the protocol)'

$ws.Range("C20").Value = 'INCL019'
$ws.Range("D20").Value = 'Inclusion Criteria'
$ws.Range("G20").Value = 'For men: agreement to remain abstinent (refrain from heterosexual intercourse) or use a condom, and agreement to refrain from donating sperm, as defined below: With a female par... (As per the protocol)'

$ws.Range("C21").Value = 'INCL020'
$ws.Range("D21").Value = 'Inclusion Criteria'
$ws.Range("G21").Value = 'Any history of prior NSCLC'

$ws.Range("C22").Value = 'EXCL001'
$ws.Range("D22").Value = 'Exclusion Criteria'
$ws.Range("G22").Value = 'Any history of prior NSCLC'

$ws.Range("C23").Value = 'EXCL002'
$ws.Range("D23").Value = 'Exclusion Criteria'
$ws.Range("G23").Value = 'NSCLC known to have a mutation in the EGFR gene or an ALK fusion oncogene are excluded from the study: – Patients with non-squamous NSCLC who have an unknown EGFR or ALK status... (As per the protocol)'

$ws.Range("C24").Value = 'EXCL003'
$ws.Range("D24").Value = 'Exclusion Criteria'
$ws.Range("G24").Value = 'If a pleural effusion is present, the following criteria must be met to exclude malignant involvement (incurable T4 disease): – When pleural fluid is visible on both the compute... (As per the protocol)'

$ws.Range("C25").Value = 'EXCL004'
$ws.Range("D25").Value = 'Exclusion Criteria'
$ws.Range("G25").Value = 'Any evidence of Stage IV disease, including, but not limited to, the following: – Pleural effusion – Pericardial effusion – Brain metastases – No history of intracranial hemorrh... (As per the protocol)'

$ws.Range("C26").Value = 'EXCL005'
$ws.Range("D26").Value = 'Exclusion Criteria'
$ws.Range("G26").Value = 'Pulmonary lymphoepithelioma-like carcinoma subtype of NSCLC'

$ws.Range("C27").Value = 'EXCL006'
$ws.Range("D27").Value = 'Exclusion Criteria'
$ws.Range("G27").Value = 'History of leptomeningeal disease'

$ws.Range("C28").Value = 'EXCL007'
$ws.Range("D28").Value = 'Exclusion Criteria'
$ws.Range("G28").Value = 'Concurrent enrollment in another clinical study, unless it is an observational (non- interventional) clinical study or the follow-up period of an interventional study'

$ws.Range("C29").Value = 'EXCL008'
$ws.Range("D29").Value = 'Exclusion Criteria'
$ws.Range("G29").Value = 'Treatment with sequential CRT for locally advanced NSCLC'

$ws.Range("C30").Value = 'EXCL009'
$ws.Range("D30").Value = 'Exclusion Criteria'
$ws.Range("G30").Value = 'Patients with locally advanced NSCLC who have progressed during or after the definite concurrent CRT prior to randomization'

$ws.Range("C31").Value = 'EXCL010'
$ws.Range("D31").Value = 'Exclusion Criteria'
$ws.Range("G31").Value = 'Any Grade > 2 unresolved toxicity from previous CRT Patients with an irreversible toxicity that is managed and is not expected to be exacerbated by study drug treatment may be i... (As per the protocol)'

$ws.Range("C32").Value = 'EXCL011'
$ws.Range("D32").Value = 'Exclusion Criteria'
$ws.Range("G32").Value = 'Grade ≥ 2 pneumonitis from prior CRT'

$ws.Range("C33").Value = 'EXCL012'
$ws.Range("D33").Value = 'Exclusion Criteria'
$ws.Range("G33").Value = 'Any concurrent chemotherapy, immunotherapy, biologic, or hormonal therapy for cancer Note: Local treatment of isolated lesions, excluding target lesions, with palliative intent... (As per the protocol)'

$ws.Range("C34").Value = 'EXCL013'
$ws.Range("D34").Value = 'Exclusion Criteria'
$ws.Range("G34").Value = 'Uncontrolled or symptomatic hypercalcemia (ionized calcium > 1.5 mmol/L, calcium > 12 mg/dL, or corrected calcium greater than ULN)'

$ws.Range("C35").Value = 'EXCL014'
$ws.Range("D35").Value = 'Exclusion Criteria'
$ws.Range("G35").Value = 'Active or history of autoimmune disease or immune deficiency, including, but not limited to, myasthenia gravis, myositis, autoimmune hepatitis, systemic lupus erythematosus, rhe... (As per the protocol)'

$ws.Range("C36").Value = 'EXCL015'
$ws.Range("D36").Value = 'Exclusion Criteria'
$ws.Range("G36").Value = 'History of idiopathic pulmonary fibrosis, organizing pneumonia (e.g., bronchiolitis obliterans), drug-induced pneumonitis, or idiopathic pneumonitis, or evidence of active pneum... (As per the protocol)'

$ws.Range("C37").Value = 'EXCL016'
$ws.Range("D37").Value = 'Exclusion Criteria'
$ws.Range("G37").Value = 'Active tuberculosis'

$ws.Range("C38").Value = 'EXCL017'
$ws.Range("D38").Value = 'Exclusion Criteria'
$ws.Range("G38").Value = 'Known clinically significant liver disease, including active viral, alcoholic, or other hepatitis, cirrhosis, and inherited liver disease, or current alcohol abuse'

$ws.Range("C39").Value = 'EXCL018'
$ws.Range("D39").Value = 'Exclusion Criteria'
$ws.Range("G39").Value = 'Significant cardiovascular disease (such as New York Heart Association Class II or greater cardiac disease, myocardial infarction, or cerebrovascular accident) within 3 months p... (As per the protocol)'

$ws.Range("C40").Value = 'EXCL019'
$ws.Range("D40").Value = 'Exclusion Criteria'
$ws.Range("G40").Value = 'Major surgical procedure, other than for diagnosis, within 4 weeks prior to initiation of study treatment, or anticipation of need for a major surgical procedure during the study'

$ws.Range("C41").Value = 'EXCL020'
$ws.Range("D41").Value = 'Exclusion Criteria'
$ws.Range("G41").Value = 'History of malignancy other than NSCLC within 5 years prior to screening, with the exception of malignancies with a negligible risk of metastasis or death (e.g., 5-year OS rate... (As per the protocol)'

$ws.Range("C42").Value = 'EXCL021'
$ws.Range("D42").Value = 'Exclusion Criteria'
$ws.Range("G42").Value = 'Severe infection within 4 weeks prior to initiation of study treatment, including, but not limited to, hospitalization for complications of infection, bacteremia, or severe pneumonia'

$ws.Range("C43").Value = 'EXCL022'
$ws.Range("D43").Value = 'Exclusion Criteria'
$ws.Range("G43").Value = 'Treatment with therapeutic oral or IV antibiotics within 2 weeks prior to initiation of study treatment Patients receiving prophylactic antibiotics (e.g., to prevent a urinary t... (As per the protocol)'

$ws.Range("C44").Value = 'EXCL023'
$ws.Range("D44").Value = 'Exclusion Criteria'
$ws.Range("G44").Value = 'Prior allogeneic stem cell or solid organ transplantation'

$ws.Range("C45").Value = 'EXCL024'
$ws.Range("D45").Value = 'Exclusion Criteria'
$ws.Range("G45").Value = 'Any other disease, metabolic dysfunction, physical examination finding, or clinical laboratory finding that contraindicates the use of an investigational drug, may affect the in... (As per the protocol)'

$ws.Range("C46").Value = 'EXCL025'
$ws.Range("D46").Value = 'Exclusion Criteria'
$ws.Range("G46").Value = 'Treatment with a live, attenuated vaccine within 4 weeks prior to initiation of study treatment, or anticipation of need for such a vaccine during study treatment or within 5 mo... (As per the protocol)'

$ws.Range("C47").Value = 'EXCL026'
$ws.Range("D47").Value = 'Exclusion Criteria'
$ws.Range("G47").Value = 'Current treatment with anti-viral therapy for HBV or HCV'

$ws.Range("C48").Value = 'EXCL027'
$ws.Range("D48").Value = 'Exclusion Criteria'
$ws.Range("G48").Value = 'Active EBV infection or known or suspected chronic active EBV infection at screening Patients with a positive EBV viral capsid antigen (VCA) IgM test at screening are excluded.... (As per the protocol)'

$ws.Range("C49").Value = 'EXCL028'
$ws.Range("D49").Value = 'Exclusion Criteria'
$ws.Range("G49").Value = 'Treatment with investigational therapy within 28 days prior to initiation of study treatment'

$ws.Range("C50").Value = 'EXCL029'
$ws.Range("D50").Value = 'Exclusion Criteria'
$ws.Range("G50").Value = 'Prior treatment with CD137 agonists or immune checkpoint blockade therapies, including anti−cytotoxic T lymphocyte−associated protein 4, anti-TIGIT, anti−PD-1, and anti−PD-L1 therapeutic antibodies'

$ws.Range("C51").Value = 'EXCL030'
$ws.Range("D51").Value = 'Exclusion Criteria'
$ws.Range("G51").Value = 'Any prior Grade ≥ 3 immune-mediated adverse event or any unresolved Grade > 1 immune-mediated adverse event while receiving any previous immunotherapy agent other than immune ch... (As per the protocol)'

$ws.Range("C52").Value = 'EXCL031'
$ws.Range("D52").Value = 'Exclusion Criteria'
$ws.Range("G52").Value = 'Treatment with systemic immunostimulatory agents (including, but not limited to, IFN and interleukin-2 [IL-2]) within 4 weeks or 5 drug- elimination half-lives (whichever is lon... (As per the protocol)'

$ws.Range("C53").Value = 'EXCL032'
$ws.Range("D53").Value = 'Exclusion Criteria'
$ws.Range("G53").Value = 'Treatment with systemic immunosuppressive medication (including, but not limited to, corticosteroids, cyclophosphamide, azathioprine, methotrexate, thalidomide, and anti−tumor n... (As per the protocol)'

$ws.Range("C54").Value = 'EXCL033'
$ws.Range("D54").Value = 'Exclusion Criteria'
$ws.Range("G54").Value = 'History of severe allergic anaphylactic reactions to chimeric or humanized antibodies or fusion proteins'

$ws.Range("C55").Value = 'EXCL034'
$ws.Range("D55").Value = 'Exclusion Criteria'
$ws.Range("G55").Value = 'Known hypersensitivity to CHO cell products or to any component of the tiragolumab or atezolizumab or durvalumab formulation'

$ws.Range("C56").Value = 'EXCL035'
$ws.Range("D56").Value = 'Exclusion Criteria'
$ws.Range("G56").Value = 'Pregnancy or breastfeeding, or intention of becoming pregnant during study treatment or within 5 months after the final dose of study treatment Women of childbearing potential m... (As per the protocol)'

$ws.Range("C57").Value = 'EXCL036'
$ws.Range("D57").Value = 'Exclusion Criteria'
$ws.Range("G57").Value = 'Any condition that, in the opinion of the investigator, would interfere with the evaluation of the study drug or interpretation of patient'
